$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.246696352958679
$ws.Range("B1").Value = 2.240676164627075
$ws.Range("C1").Value = 2.884564161300659
$ws.Range("D1").Value = 3.342852830886841
$ws.Range("E1").Value = 1.931156635284424
